$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = -19.86529999999998
$ws.Range("B4").Value = 8.878999999999994
$ws.Range("A6").Value = -22.54830000000001
$ws.Range("A7").Value = -19.4337
$ws.Range("D7").Value = -7.742199999999997
$ws.Range("D8").Value = -8.453800000000003
$ws.Range("B9").Value = 5.926700000000003
$ws.Range("D10").Value = -7.404800000000002
$ws.Range("B12").Value = 4.889899999999997
$ws.Range("D13").Value = -8.492799999999995
$ws.Range("A16").Value = -21.85410000000001
$ws.Range("D16").Value = -8.645500000000004
$ws.Range("B17").Value = 6.112099999999994
$ws.Range("B18").Value = 6.706699999999994
$ws.Range("B19").Value = 9.408599999999989
$ws.Range("A20").Value = -22.71090000000002
$ws.Range("B20").Value = 5.408899999999999
$ws.Range("B26").Value = 5.380600000000007
$ws.Range("A28").Value = -21.94589999999999
$ws.Range("A29").Value = -21.26429999999996
$ws.Range("D30").Value = -7.087499999999993
$ws.Range("B31").Value = 4.840599999999998
$ws.Range("A32").Value = -21.18160000000001
$ws.Range("B39").Value = 9.441800000000002
$ws.Range("A40").Value = -20.02439999999998
$ws.Range("B40").Value = 8.745899999999995
$ws.Range("D40").Value = -8.075800000000005
$ws.Range("B41").Value = 9.533799999999989
$ws.Range("B42").Value = 9.535899999999987
$ws.Range("B43").Value = 6.171800000000002
$ws.Range("D44").Value = -6.682100000000005
$ws.Range("A46").Value = -21.7272
$ws.Range("B47").Value = 4.905599999999996
$ws.Range("B48").Value = 5.747800000000002
$ws.Range("A51").Value = -21.68569999999998
$ws.Range("A52").Value = -22.29629999999999
$ws.Range("A57").Value = -21.92440000000002
$ws.Range("A59").Value = -22.08100000000001
$ws.Range("A62").Value = -21.93130000000001
$ws.Range("B63").Value = 4.833099999999998
$ws.Range("B64").Value = 5.344600000000002
$ws.Range("A66").Value = -21.37909999999998
$ws.Range("A73").Value = -20.20539999999998
$ws.Range("A74").Value = -21.70329999999998
$ws.Range("B76").Value = 6.029199999999998
$ws.Range("B81").Value = 5.002300000000004
$ws.Range("B89").Value = 5.340599999999998
$ws.Range("D89").Value = -8.642099999999999
$ws.Range("D91").Value = -7.918999999999998
$ws.Range("A92").Value = -21.63840000000001
$ws.Range("B94").Value = 5.032199999999991
$ws.Range("A100").Value = -22.00130000000001